$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - "time_taken", styled like other headers (same style as E1)
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "time_taken"

# Data cells F2:F42 - time_taken values as text strings
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:41:00.828146"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:41:00.828155"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:41:00.828158"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:41:00.828160"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:41:00.828163"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:41:00.828165"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:41:00.828167"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:41:00.828170"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:41:00.828172"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:41:00.828174"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:41:00.828176"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:41:00.828179"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:41:00.828181"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:41:00.828183"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:41:00.828185"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:41:00.828187"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:41:00.828190"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:41:00.828192"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:41:00.828194"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:41:00.828196"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:41:00.828198"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:41:00.828201"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:41:00.828203"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:41:00.828205"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:41:00.828208"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:41:00.828210"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:41:00.828212"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:41:00.828215"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:41:00.828217"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:41:00.828219"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:41:00.828222"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:41:00.828224"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:41:00.828227"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:41:00.828229"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:41:00.828232"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:41:00.828234"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:41:00.828236"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:41:00.828239"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:41:00.828243"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:41:00.828246"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:41:00.828248"

$wb.Save()
